# functions in spark added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "windo_open" -> "window_open" in B44
$ws.Range("B44").Value = "window_open"

# Remove the helper formula from B27 (was =+B2, showing "start_drive")
$ws.Range("B27").ClearContents()

# Leave the selection on B27, matching the recorded cursor position
$ws.Range("B27").Select()
